# Updated cryptos list on Mon Sep 11 02:05:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates, keyed by row number. These are stored as text
# in the sheet (e.g. "25.855.67", "0.0₃0735"), so force a text number
# format before writing the value to stop Excel from auto-converting
# numeric-looking strings into real numbers. Reset the style afterwards so
# no stray cell style/format is left behind.
$priceUpdates = @{
    2  = "25.855.67"
    3  = "1.618.40"
    5  = "213.17"
    6  = "0.498"
    9  = "0.0614"
    10 = "18.34"
    12 = "1.843.88"
    13 = "1.612.53"
    14 = "4.12"
    16 = "25.868.19"
    17 = "61.28"
    18 = "0.0₃0735"
    20 = "190.72"
    22 = "9.45"
    25 = "143.44"
    27 = "1.71"
    29 = "15.17"
    31 = "0.0476"
    34 = "2.41"
    36 = "1.120.49"
    37 = "0.838"
    39 = "0.0153"
    40 = "0.509"
    41 = "97.87"
    42 = "1.754.61"
    43 = "0.747"
    44 = "5.06"
    45 = "0.0₆0112"
    46 = "53.98"
    48 = "0.0520"
    49 = "0.412"
    51 = "7.44"
}

# Column E (Volume 1h) updates, keyed by row number. Values already carry
# leading/trailing spaces and a "%" sign, so Excel naturally stores them
# as text -- no special handling required.
$volumeUpdates = @{
    2  = "  -0.10%  "
    3  = "  -0.89%  "
    4  = "  +0.62%  "
    5  = "  -0.48%  "
    6  = "  -1.22%  "
    7  = "  +0.64%  "
    8  = "  -1.05%  "
    9  = "  -3.07%  "
    10 = "  -5.18%  "
    11 = "  -0.20%  "
    12 = "  -0.83%  "
    13 = "  -1.17%  "
    14 = "  -2.36%  "
    16 = "  -0.14%  "
    17 = "  -2.48%  "
    18 = "  -2.75%  "
    19 = "  +0.51%  "
    20 = "  -1.11%  "
    21 = "  -1.59%  "
    22 = "  -2.43%  "
    23 = "  -1.63%  "
    24 = "  +1.52%  "
    25 = "  -0.07%  "
    26 = "  +0.65%  "
    27 = "  -2.81%  "
    28 = "  -2.05%  "
    29 = "  -1.61%  "
    30 = "  -0.85%  "
    31 = "  -2.83%  "
    32 = "  -3.83%  "
    33 = "  -4.73%  "
    34 = "  -1.37%  "
    35 = "  -2.84%  "
    36 = "  -0.08%  "
    37 = "  -6.28%  "
    38 = "  -2.36%  "
    39 = "  -1.59%  "
    40 = "  -4.35%  "
    41 = "  -0.24%  "
    42 = "  -0.68%  "
    43 = "  -5.77%  "
    44 = "  -4.92%  "
    45 = "  -1.61%  "
    46 = "  -3.46%  "
    47 = "  +0.25%  "
    48 = "  +0.22%  "
    49 = "  +0.11%  "
    50 = "  +0.82%  "
    51 = "  -3.10%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

# Rows 17 and 18 swap coin identity (Coin name + Link); the A (rank) column
# stays put. Price/Volume for these rows are already handled above.
$ws.Cells.Item(17, 2).Value = "Litecoin"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"

$ws.Cells.Item(18, 2).Value = "ShibaInu"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"

# Rows 39 and 40 swap coin identity (Coin name + Link) as well.
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

$ws.Cells.Item(40, 2).Value = "ImmutableX"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
